# "cambio con todas las variables sin depuracion"
# Replace the predicted ("y_pred") values in column B with the new,
# un-filtered prediction run. Column A (dates) and the headers
# (A1="fa", B1="y_pred") are unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(
    -34.85344466089452,
    989.2614185798475,
    199.3866436327563,
    -87.17837029581527,
    169.5373767099568,
    634.0430841702744,
    1122.774828374219
)

$row = 2
foreach ($v in $newValues) {
    $ws.Cells.Item($row, 2).Value = $v
    $row++
}
